$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell references to their new text values. Using a literal-text
# write (NumberFormat "@" while assigning, then restoring the "Normal"
# cell style) ensures values like "531.27" or "61.039.33" are stored as
# text, matching the source data which is not numeric (dotted price
# strings, percentage strings, etc.), instead of being auto-coerced into
# numbers by Excel's usual text-looks-like-a-number parsing.
$updates = @{
    "D2" = "61.039.33"
    "E2" = "  +0.30%  "
    "D3" = "2.672.07"
    "E3" = "  +1.98%  "
    "E4" = "  -0.01%  "
    "D5" = "531.27"
    "E5" = "  +3.37%  "
    "D6" = "155.90"
    "E6" = "  +0.40%  "
    "D7" = "0.997"
    "E7" = "  -0.01%  "
    "D8" = "0.586"
    "E8" = "  -0.19%  "
    "D9" = "6.58"
    "E9" = "  -3.35%  "
    "E10" = "  +3.91%  "
    "D11" = "0.354"
    "E11" = "  +1.89%  "
    "E12" = "  -0.43%  "
    "D13" = "3.141.58"
    "E13" = "  +2.00%  "
    "D14" = "61.035.71"
    "E14" = "  +0.38%  "
    "D15" = "22.11"
    "E15" = "  +1.64%  "
    "E16" = "  +1.39%  "
    "D17" = "2.685.65"
    "E17" = "  +2.16%  "
    "E18" = "  +0.96%  "
    "D19" = "355.57"
    "E19" = "  -0.13%  "
    "D20" = "10.73"
    "E20" = "  +0.95%  "
    "D21" = "6.35"
    "E21" = "  +2.45%  "
    "D22" = "0.998"
    "E22" = "  -0.18%  "
    "D23" = "61.80"
    "E23" = "  +1.40%  "
    "D25" = "0.168"
    "E25" = "  +1.02%  "
    "D26" = "0.999"
    "E26" = "  +0.28%  "
    "D27" = "0.0₃0862"
    "E27" = "  +1.62%  "
    "D28" = "7.39"
    "E28" = "  +0.09%  "
    "E29" = "  -0.02%  "
    "E30" = "  +4.20%  "
    "D31" = "19.59"
    "E31" = "  +0.65%  "
    "E32" = "  +2.98%  "
    "D33" = "149.93"
    "E33" = "  -1.56%  "
    "D34" = "4.15"
    "E34" = "  +3.35%  "
    "E35" = "  +0.73%  "
    "E36" = "  +7.97%  "
    "E37" = "  +1.77%  "
    "E38" = "  +0.12%  "
    "D39" = "36.87"
    "E39" = "  +1.38%  "
    "D40" = "306.44"
    "E40" = "  +4.09%  "
    "D41" = "3.79"
    "E41" = "  +0.45%  "
    "D42" = "0.651"
    "E42" = "  +4.05%  "
    "B43" = "Stellar"
    "C43" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D43" = "0.102"
    "E43" = "  +0.51%  "
    "B44" = "EnergySwap"
    "C44" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D44" = "20.60"
    "E44" = "  +3.48%  "
    "D45" = "0.0567"
    "E45" = "  +1.73%  "
    "E46" = "  +0.06%  "
    "E47" = "  +2.85%  "
    "D48" = "4.92"
    "E48" = "  -0.67%  "
    "D49" = "19.17"
    "E49" = "  +7.61%  "
    "E50" = "  +0.41%  "
    "D51" = "2.003.82"
    "E51" = "  +0.00%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
